# Update "想去人数" (want-to-go count) values on both the "展览" sheet and
# the "全部类型" sheet, which each carry a duplicated copy of the same rows.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 3248
$wsExhibit.Range("F7").Value = 4735
$wsExhibit.Range("F21").Value = 4694
$wsExhibit.Range("F27").Value = 1181
$wsExhibit.Range("F37").Value = 757

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value = 3248
$wsAll.Range("F11").Value = 4735
$wsAll.Range("F26").Value = 4694
$wsAll.Range("F32").Value = 1181
$wsAll.Range("F43").Value = 757
